# Windchill TestDataInput.xlsx update
# - Adds three new worksheets (BaseLine_Creation, ChangeNotice, Product_Creation)
#   between "Credentials" and "Part_Creation".
# - Adds a new credentials row (testuser2 / 123) to the "Credentials" sheet.

$wb = $excel.ActiveWorkbook
$credentials = $wb.Worksheets.Item("Credentials")

# --- New sheet: BaseLine_Creation -------------------------------------------------
$baseline = $wb.Worksheets.Add($null, $credentials)
$baseline.Name = "BaseLine_Creation"

$baseline.Range("A1").Value = "Name"
$baseline.Range("B1").Value = "Description"
$baseline.Range("A2").Value = "BaseLine_1"
$baseline.Range("B2").Value = "For testing"

$baseline.Columns.Item(1).ColumnWidth = 15.608072916666666
$baseline.Columns.Item(2).ColumnWidth = 19.944010416666668

$baseline.Range("B2").Select()

# --- New sheet: ChangeNotice --------------------------------------------------------
$changeNotice = $wb.Worksheets.Add($null, $baseline)
$changeNotice.Name = "ChangeNotice"

$changeNotice.Range("A1").Value = "CN_Name"
$changeNotice.Range("B1").Value = "CT_Name"
$changeNotice.Range("C1").Value = "Approver"
$changeNotice.Range("D1").Value = "Reviewer"
$changeNotice.Range("A2").Value = "Change_Notice1"
$changeNotice.Range("B2").Value = "Change_task1"
$changeNotice.Range("C2").Value = "testuser2"
$changeNotice.Range("D2").Value = "testuser2"

$changeNotice.Columns.Item(1).ColumnWidth = 14.053385416666666
$changeNotice.Columns.Item(2).ColumnWidth = 17.944010416666668
$changeNotice.Columns.Item(3).ColumnWidth = 20.498697916666668

$changeNotice.Range("D2").Select()

# --- New sheet: Product_Creation ----------------------------------------------------
$productCreation = $wb.Worksheets.Add($null, $changeNotice)
$productCreation.Name = "Product_Creation"

$productCreation.Range("A1").Value = "Name"
$productCreation.Range("B1").Value = "Description"
$productCreation.Range("A2").Value = "Product_1"
$productCreation.Range("B2").Value = "For testing"

$productCreation.Columns.Item(1).ColumnWidth = 12.830729166666666
$productCreation.Columns.Item(2).ColumnWidth = 18.608072916666668

$productCreation.Range("B2").Select()

# --- Credentials: append a new user row ---------------------------------------------
$credentials.Range("A4").Value = "testuser2"
$credentials.Range("B4").Value = 123

$credentials.Activate()
$credentials.Range("A4").Select()
